# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement "Periodo Mora" / "Valor Mora" table (rows 16-22) is
# refreshed: the old periods are removed and the new ones added, which in
# effect reverses the order of the seven periods and moves the 38000 "Valor
# Mora" amount from period 2005 (last row) to period 1911 (first row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2005", "2004", "2003", "2002", "2001", "1912", "1911")
$valores = @(38000, 60000, 60000, 60000, 60000, 60000, 60000)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}

# Column widths were re-fit (best-fit) against the refreshed data/fonts.
$ws.Columns.Item(2).ColumnWidth  = 17.709635416666668
$ws.Columns.Item(3).ColumnWidth  = 15.893229166666666
$ws.Columns.Item(5).ColumnWidth  = 12.709635416666666
$ws.Columns.Item(6).ColumnWidth  = 9.346354166666666
$ws.Columns.Item(7).ColumnWidth  = 13.529947916666666
$ws.Columns.Item(8).ColumnWidth  = 18.529947916666668
$ws.Columns.Item(9).ColumnWidth  = 17.256510416666668
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666

$wb.Save()
